$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 3090.25
$ws.Range("I29").Value = 1944.4286
$ws.Range("K29").Value = 5833.2858
$ws.Range("M29").Value = -5552.2858

$ws.Range("H76").Value = 3706837
$ws.Range("I76").Value = 3971332.5
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 3971332.5
$ws.Range("L76").Value = 3900
$ws.Range("M76").Value = -3971017.5
$ws.Range("N76").Value = -4530

$ws.Range("H79").Value = 3706837
$ws.Range("I79").Value = 3971332.5
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 3971332.5
$ws.Range("L79").Value = 3900
$ws.Range("M79").Value = -3970240.5
$ws.Range("N79").Value = -6084

$ws.Range("H116").Value = 10649282
$ws.Range("I116").Value = 23063148
$ws.Range("J116").Value = 8825.714
$ws.Range("K116").Value = 23063148
$ws.Range("L116").Value = 8825.714
$ws.Range("M116").Value = -23059706
$ws.Range("N116").Value = -15709.714

$ws.Range("H137").Value = 1886.7059
$ws.Range("I137").Value = 1969.1428
$ws.Range("K137").Value = 5907.428400000001
$ws.Range("M137").Value = -3357.428400000001

$ws.Range("H138").Value = 1727.3776
$ws.Range("J138").Value = 1982.1948
$ws.Range("L138").Value = 5946.5844
$ws.Range("N138").Value = -16226.5844

$ws.Range("H141").Value = 1750.6818
$ws.Range("I141").Value = 1595.9524
$ws.Range("K141").Value = 4787.857199999999
$ws.Range("M141").Value = 392.1428000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 15179.556
$ws.Range("I74").Value = 2352
$ws.Range("J74").Value = 40834.668
$ws.Range("K74").Value = 2352
$ws.Range("L74").Value = 40834.668
$ws.Range("M74").Value = -1478
$ws.Range("N74").Value = -42582.668

$ws.Range("H77").Value = 15179.556
$ws.Range("I77").Value = 2352
$ws.Range("J77").Value = 40834.668
$ws.Range("K77").Value = 11760
$ws.Range("L77").Value = 204173.34
$ws.Range("M77").Value = -7392
$ws.Range("N77").Value = -212909.34

$ws.Range("H122").Value = 6797.0557
$ws.Range("I122").Value = 7643
$ws.Range("K122").Value = 22929
$ws.Range("M122").Value = -20479

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 34191
$ws.Range("J81").Value = 34191
$ws.Range("L81").Value = 34191
$ws.Range("N81").Value = -36313

$ws.Range("H84").Value = 34191
$ws.Range("J84").Value = 34191
$ws.Range("L84").Value = 102573
$ws.Range("N84").Value = -113181

$ws.Range("H134").Value = 3231.8948
$ws.Range("I134").Value = 2363.4666
$ws.Range("J134").Value = 6488.5
$ws.Range("K134").Value = 7090.399800000001
$ws.Range("L134").Value = 19465.5
$ws.Range("M134").Value = -4555.399800000001
$ws.Range("N134").Value = -24535.5

$ws.Range("H135").Value = 88833
$ws.Range("J135").Value = 88833
$ws.Range("L135").Value = 88833
$ws.Range("N135").Value = -98973

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3727.6738
$ws.Range("I31").Value = 1585.625
$ws.Range("J31").Value = 6064.4546
$ws.Range("K31").Value = 1585.625
$ws.Range("L31").Value = 6064.4546
$ws.Range("M31").Value = -1290.625
$ws.Range("N31").Value = -6654.4546

$ws.Range("H34").Value = 3727.6738
$ws.Range("I34").Value = 1585.625
$ws.Range("J34").Value = 6064.4546
$ws.Range("K34").Value = 1585.625
$ws.Range("L34").Value = 6064.4546
$ws.Range("M34").Value = -1383.625
$ws.Range("N34").Value = -6468.4546

$ws.Range("H122").Value = 1588
$ws.Range("I122").Value = 1222
$ws.Range("J122").Value = 2228.5
$ws.Range("K122").Value = 3666
$ws.Range("L122").Value = 6685.5
$ws.Range("M122").Value = -1216
$ws.Range("N122").Value = -11585.5

$ws.Range("H123").Value = 32000
$ws.Range("J123").Value = 32000
$ws.Range("L123").Value = 32000
$ws.Range("N123").Value = -41800

$ws.Range("H124").Value = 44081
$ws.Range("J124").Value = 44081
$ws.Range("L124").Value = 44081
$ws.Range("N124").Value = -48991

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2681.5264
$ws.Range("I5").Value = 1972
$ws.Range("J5").Value = 3320.1
$ws.Range("K5").Value = 5916
$ws.Range("L5").Value = 9960.299999999999
$ws.Range("M5").Value = -5804
$ws.Range("N5").Value = -10184.3

$ws.Range("H122").Value = 1088
$ws.Range("I122").Value = 632.5
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5692.5
$ws.Range("L122").Value = 17991
$ws.Range("M122").Value = -3242.5
$ws.Range("N122").Value = -22891

$ws.Range("H131").Value = 2925.4922
$ws.Range("J131").Value = 3052.8547
$ws.Range("L131").Value = 9158.5641
$ws.Range("N131").Value = -19238.5641

$ws.Range("H135").Value = 2681.5264
$ws.Range("I135").Value = 1972
$ws.Range("J135").Value = 3320.1
$ws.Range("K135").Value = 17748
$ws.Range("L135").Value = 29880.9
$ws.Range("M135").Value = -15213
$ws.Range("N135").Value = -34950.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3184
$ws.Range("I7").Value = 2577.7778
$ws.Range("J7").Value = 3525
$ws.Range("K7").Value = 2577.7778
$ws.Range("L7").Value = 3525
$ws.Range("M7").Value = -2465.7778
$ws.Range("N7").Value = -3749

$ws.Range("H22").Value = 14552.857
$ws.Range("I22").Value = 160
$ws.Range("J22").Value = 25347.5
$ws.Range("K22").Value = 160
$ws.Range("L22").Value = 25347.5
$ws.Range("M22").Value = 135
$ws.Range("N22").Value = -25937.5

$ws.Range("H27").Value = 14552.857
$ws.Range("I27").Value = 160
$ws.Range("J27").Value = 25347.5
$ws.Range("K27").Value = 160
$ws.Range("L27").Value = 25347.5
$ws.Range("M27").Value = -53
$ws.Range("N27").Value = -25561.5

$ws.Range("H40").Value = 3450
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3473.6843
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3473.6843
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -3745.6843

$ws.Range("H126").Value = 3184
$ws.Range("I126").Value = 2577.7778
$ws.Range("J126").Value = 3525
$ws.Range("K126").Value = 7733.3334
$ws.Range("L126").Value = 10575
$ws.Range("M126").Value = -5263.3334
$ws.Range("N126").Value = -15515

$ws.Range("H132").Value = 3581.2068
$ws.Range("I132").Value = 2091.0625
$ws.Range("K132").Value = 6273.1875
$ws.Range("M132").Value = -3743.1875

$ws.Range("H133").Value = 42257
$ws.Range("J133").Value = 42257
$ws.Range("L133").Value = 42257
$ws.Range("N133").Value = -47317

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 17375
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 17375
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 17375
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -17855

$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H126").Value = 48767.477
$ws.Range("I126").Value = 84234.914
$ws.Range("K126").Value = 252704.742
$ws.Range("M126").Value = -250234.742
